# Update countries & provincias Spain
# Applies the data refresh to the "paises" sheet:
#  - Updates case counts for several countries (new day's data)
#  - Because some countries' totals overtook their neighbours, those rows
#    now show a different country name (the ranking swapped rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name changes (rank reshuffles) -------------------------------
$ws.Range("A28").Value = "Chile"
$ws.Range("A29").Value = "Singapur"

$ws.Range("A74").Value = "Camerun"
$ws.Range("A75").Value = "Azerbaiyan"
$ws.Range("A76").Value = "Bosnia y Herzegovina"

$ws.Range("A193").Value = "San Vicente y las Granadinas"
$ws.Range("A194").Value = "Namibia"

$ws.Range("A217").Value = "San Pedro y Miquelon"
$ws.Range("A218").Value = "Comoras"

# --- Numeric data updates (Casos totales, Nuevos casos, Casos activos, ---
# --- Recuperados, Casos criticos, Muertes hoy, Muertes) -------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1135984
$ws.Range("C4").Value = 4954
$ws.Range("E4").Value = 907943
$ws.Range("G4").Value = 506
$ws.Range("H4").Value = 66259

# Row 9: Alemania
$ws.Range("B9").Value = 164380
$ws.Range("C9").Value = 303
$ws.Range("E9").Value = 28644

# Row 19: India
$ws.Range("B19").Value = 37776
$ws.Range("C19").Value = 519
$ws.Range("E19").Value = 26546

# Row 28: now Chile
$ws.Range("B28").Value = 18435
$ws.Range("C28").Value = 1427
$ws.Range("D28").Value = 9572
$ws.Range("E28").Value = 8616
$ws.Range("F28").Value = 386
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = 247

# Row 29: now Singapur
$ws.Range("B29").Value = 17548
$ws.Range("C29").Value = 447
$ws.Range("D29").Value = 1347
$ws.Range("E29").Value = 16184
$ws.Range("F29").Value = 24
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 17

# Row 46: Chequia
$ws.Range("B46").Value = 7750
$ws.Range("C46").Value = 13
$ws.Range("D46").Value = 3446
$ws.Range("E46").Value = 4059
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 245

# Row 47: Republica Dominicana
$ws.Range("B47").Value = 7578
$ws.Range("C47").Value = 290
$ws.Range("D47").Value = 1481
$ws.Range("E47").Value = 5771
$ws.Range("G47").Value = 13
$ws.Range("H47").Value = 326

# Row 74: now Camerun
$ws.Range("B74").Value = 2077
$ws.Range("C74").Value = 245
$ws.Range("D74").Value = 953
$ws.Range("E74").Value = 1060
$ws.Range("F74").Value = 12
$ws.Range("G74").Value = 3
$ws.Range("H74").Value = 64

# Row 75: now Azerbaiyan
$ws.Range("B75").Value = 1894
$ws.Range("C75").Value = 40
$ws.Range("D75").Value = 1411
$ws.Range("E75").Value = 458
$ws.Range("F75").Value = 17
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 25

# Row 76: now Bosnia y Herzegovina
$ws.Range("B76").Value = 1839
$ws.Range("C76").Value = 58
$ws.Range("D76").Value = 779
$ws.Range("E76").Value = 988
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 72

# Row 117: Jordania
$ws.Range("B117").Value = 460
$ws.Range("C117").Value = 1
$ws.Range("D117").Value = 367
$ws.Range("E117").Value = 84
